$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append 8 new rows (50-57) for subject S07 / anderson, session 2013-07-05,
# runs 1-8 -- mirroring the existing per-subject blocks already in the sheet.
# ---------------------------------------------------------------------------

# First, stamp each new row with the same formatting (styles) as an existing,
# fully-styled data row (row 10: centered text/date styles on A,B,C,F).
for ($i = 0; $i -lt 8; $i++) {
  $r = 50 + $i
  $ws.Range("A10:F10").Copy($ws.Range("A" + $r + ":F" + $r))
}

$fileNames = @(
  "2013-07-05-14-56-29-run1",
  "2013-07-05-15-03-20-run2",
  "2013-07-05-15-09-42-run3",
  "2013-07-05-15-21-37-run4",
  "2013-07-05-15-36-23-run5",
  "2013-07-05-15-42-28-run6",
  "2013-07-05-15-48-04-run7",
  "2013-07-05-15-54-43-run8"
)

# Column D (sessionDirectory) - identical for every new row.
for ($i = 0; $i -lt 8; $i++) {
  $r = 50 + $i
  $ws.Cells.Item($r, 4).Value = "2013-07-05-anderson"
}

# Column A (subjectTag).
for ($i = 0; $i -lt 8; $i++) {
  $r = 50 + $i
  $ws.Cells.Item($r, 1).Value = "S07"
}

# Column B (subjectName).
for ($i = 0; $i -lt 8; $i++) {
  $r = 50 + $i
  $ws.Cells.Item($r, 2).Value = "anderson"
}

# Column E (fileName) - one distinct value per row.
for ($i = 0; $i -lt 8; $i++) {
  $r = 50 + $i
  $ws.Cells.Item($r, 5).Value = $fileNames[$i]
}

# Column C (date) and Column F (run number).
for ($i = 0; $i -lt 8; $i++) {
  $r = 50 + $i
  $ws.Cells.Item($r, 3).Value = 41401
  $ws.Cells.Item($r, 6).Value = $i + 1
}

# ---------------------------------------------------------------------------
# View tweaks: selection moves to L34, and column D widens to fit the new,
# longer "sessionDirectory" text.
# ---------------------------------------------------------------------------
$ws.Range("L34").Select()
$ws.Columns("D:D").ColumnWidth = 18.8333333333333
